$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the paragraph whose Range.Text contains $needle (first match)
# and return its 1-based index within $d.Paragraphs.
# ---------------------------------------------------------------------------
function Find-ParagraphIndex($needle) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

$lq    = [char]0x201C   # “
$rq    = [char]0x201D   # ”
$apos  = [char]0x2019   # '  (right single quote / apostrophe)
$lsq   = [char]0x2018   # '  (left single quote)

# ---------------------------------------------------------------------------
# Change 1: remove the whole paragraph
#   "If she is above 40 she will say "maybe later""
# (including the proofErr-wrapped "40" run) that sits between the
# "If GIRL is not above a 50..." paragraph and the
# "Between 40-30 lowers Happiness by 5 points" paragraph.
# ---------------------------------------------------------------------------
$idx = Find-ParagraphIndex("If she is above")
if ($idx -gt 0) {
    $d.Paragraphs.Item($idx).Range.Delete()
}

# ---------------------------------------------------------------------------
# Changes 2-5: the <w:lastRenderedPageBreak/> marker that currently sits as
# the first child of the first run of paragraph "A" below needs to move to
# become the first child of the first run of the very next paragraph "B"
# (a plain rendering artifact that shifted because of upstream edits).
# We implement the move by rewriting both paragraphs (A without the break,
# B with the break prefixed to its first run) via Range.InsertXML, which
# replaces the exact range it is called on.
# ---------------------------------------------------------------------------

function Move-LastRenderedPageBreak($needleA, $xmlA, $xmlB) {
    $iA = Find-ParagraphIndex($needleA)
    if ($iA -le 0) { return }
    $pA = $d.Paragraphs.Item($iA)
    $pB = $d.Paragraphs.Item($iA + 1)
    $rng = $d.Range($pA.Range.Start, $pB.Range.End)
    $rng.InsertXML($xmlA + $xmlB)
}

# --- Pair 1: "Has trait= ... exercise too!" / "DNT = "I can tell ☺"" -------
$xmlA1 = "<w:p><w:r><w:t>Has trait= ${lq}that${apos}s awesome! I love exercise too!${rq}</w:t></w:r></w:p>"
$xmlB1 = "<w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space=`"preserve`">DNT = ${lq}I can tell </w:t></w:r><w:r><w:sym w:font=`"Wingdings`" w:char=`"F04A`"/></w:r><w:r><w:t>${rq}</w:t></w:r></w:p>"
Move-LastRenderedPageBreak "Has trait=" $xmlA1 $xmlB1

# --- Pair 2: "Passed = WOW! You're amazing" / "Failed "OH no! ..."" -------
$xmlA2 = "<w:p><w:pPr><w:ind w:left=`"360`"/></w:pPr><w:r><w:t xml:space=`"preserve`">Passed = WOW! You${apos}re amazing </w:t></w:r></w:p>"
$xmlB2 = "<w:p><w:pPr><w:ind w:left=`"360`"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Failed ${lq}OH no! are you ok? I think you ripped your pants${rq}</w:t></w:r></w:p>"
Move-LastRenderedPageBreak "Passed = WOW" $xmlA2 $xmlB2

# --- Pair 3: "Depends on my mood +1" / ""fair"" ---------------------------
$xmlA3 = "<w:p><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr><w:r><w:t>Depends on my mood +1</w:t></w:r></w:p>"
$xmlB3 = "<w:p><w:r><w:lastRenderedPageBreak/><w:t>${lq}fair${rq}</w:t></w:r></w:p>"
Move-LastRenderedPageBreak "Depends on my mood" $xmlA3 $xmlB3

# --- Pair 4: "I built a robot..." / "Passed and HT = "HAHAHAHA! ...""  ----
$xmlA4 = "<w:p><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr>" +
         "<w:r><w:t>I built</w:t></w:r>" +
         "<w:r><w:t xml:space=`"preserve`"> a</w:t></w:r>" +
         "<w:r><w:t xml:space=`"preserve`"> robot that passes butter (intelligence </w:t></w:r>" +
         "<w:r><w:t xml:space=`"preserve`">80%) +10 if has trait ${lq}likes rick and </w:t></w:r>" +
         "<w:proofErr w:type=`"spellStart`"/>" +
         "<w:r><w:t>morty</w:t></w:r>" +
         "<w:proofErr w:type=`"spellEnd`"/>" +
         "<w:r><w:t xml:space=`"preserve`">${rq} +5 if passed </w:t></w:r>" +
         "<w:r><w:t xml:space=`"preserve`">-3 </w:t></w:r>" +
         "<w:r><w:t>if failed</w:t></w:r></w:p>"
$xmlB4 = "<w:p><w:pPr><w:ind w:left=`"360`"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Passed and HT = ${lq}HAHAHAHA! ${lsq}Oh my god${apos} ${lq}</w:t></w:r></w:p>"
Move-LastRenderedPageBreak "I built" $xmlA4 $xmlB4
